# "feat: update draft exploration"
#
# The sheet has a header row (A1:K1) with an autoFilter already defined
# over A1:K110. Column K ("interesting", the 11th column -> filter colId 10)
# is a boolean flag. This edit narrows the draft exploration down to the
# "interesting" rows by filtering column K to TRUE, which hides every row
# whose K value isn't TRUE (FALSE or blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K is the 11th column of the already-existing autoFilter range
# (A1:K110). Field 11 (1-based) = column K. Using the xlFilterValues
# operator (7) with an explicit criteria array gives a discrete-values
# filter (<filters><filter val="TRUE"/></filters>) rather than a
# custom/comparison filter, matching a manual "check just TRUE" pick from
# the AutoFilter dropdown.
$ws.Range("A1:K110").AutoFilter(11, @("TRUE"), 7)

# Scroll the view roughly into the filtered data (best effort - mirrors
# the author scrolling back up/over after re-filtering).
try {
    $win = $excel.ActiveWindow
    $win.ScrollColumn = 3
    $win.ScrollRow = 1
} catch {
    # View-scroll state isn't critical to the data edit; ignore if the
    # host doesn't support it.
}
